# This script reproduces the edit described in the commit: a new weekly
# price record (week of 2023-11-28, serial 45258) is inserted as a new
# data row right before the current row 568 ("Fruta / hortaliza, semanal").
# All the subsequent data rows (old 568..687) shift down by one row to
# 569..688, and the worksheet's used range grows from A1:R687 to A1:R688.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 568; this pushes the previous
# rows 568-687 down to 569-688 and carries formatting (e.g. the date
# number format in column D) down from the row above automatically.
$ws.Rows("568").Insert()

# Populate the newly inserted row 568 with the new record's data.
$ws.Range("A568").Value = 5
$ws.Range("B568").Value = "Macroferia Regional de Talca"
$ws.Range("C568").Value = "Maule"
$ws.Range("D568").Value = 45258
$ws.Range("E568").Value = 7
$ws.Range("F568").Value = 100112032
$ws.Range("G568").Value = "Zapallo italiano"
$ws.Range("H568").Value = "Sin especificar"
$ws.Range("I568").Value = "Primera"
$ws.Range("J568").Value = 500
$ws.Range("K568").Value = 10000
$ws.Range("L568").Value = 11000
$ws.Range("M568").Value = 10600
$ws.Range("N568").Value = "$/caja 50 unidades"
$ws.Range("O568").Value = "Región del Maule"
$ws.Range("P568").Value = 212
$ws.Range("Q568").Value = 50
$ws.Range("R568").Value = "Hortaliza"
